$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add new row 15 - a new "2023" data row, built by duplicating row 14
#    (same values + same formatting) and then updating the two cells that
#    differ (D = item index "14", E = year "2023").
# ---------------------------------------------------------------------------
$ws.Range("A14:V14").Copy()
$ws.Range("A15:V15").PasteSpecial()
$ws.Range("A14:V14").Copy()
$ws.Range("A15:V15").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# D15 / E15 need to hold the new values as text (matching how every other
# "numeric looking" value in this sheet is stored as text), so force a text
# number format before assigning, then restore the original General format
# used by the other cells in the row.
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2023"

$ws.Range("D14:E14").Copy()
$ws.Range("D15:E15").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. New API query was run - the short-url ("B" column) changes for every
#    data row, and the "oip" column ("U") null marker is now rendered as "-"
#    instead of the literal text "null" (same change also applies to the new
#    row 15 that was just created above).
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 2).Value = "9I6Fck"
    $ws.Cells.Item($r, 21).Value = "-"
}

# ---------------------------------------------------------------------------
# 3. The "hst" column ("V") used to report a bare "0" for every historic row;
#    it now reports "-" (left aligned, like the other not-applicable columns)
#    for every row except the most recent two (2022 and 2023), which keep
#    reporting "0" right aligned.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 22)
    $cell.Value = "-"
    $cell.HorizontalAlignment = -4131
}
